$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K data to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (shifted original column D) into new D:E,
# restricted to the three data blocks (Income Statement / Balance Sheet / Cash Flow)
# so header/section-title rows that never had D:K cells stay untouched.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the two new quarter columns (D = 2018-12-31, E = 2018-09-30) with reported figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 361700
$ws.Range("E8").Value = 357900
$ws.Range("D9").Value = 73800
$ws.Range("E9").Value = 70900
$ws.Range("D10").Value = 287900
$ws.Range("E10").Value = 287000
$ws.Range("D12").Value = 20300
$ws.Range("E12").Value = 20600
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 19900
$ws.Range("E15").Value = 19100
$ws.Range("D17").Value = 191900
$ws.Range("E17").Value = 181500
$ws.Range("D18").Value = 169800
$ws.Range("E18").Value = 176400
$ws.Range("D20").Value = 53400
$ws.Range("E20").Value = 6400
$ws.Range("D21").Value = 243100
$ws.Range("E21").Value = 201900
$ws.Range("D22").Value = 35900
$ws.Range("E22").Value = 35900
$ws.Range("D23").Value = 187300
$ws.Range("E23").Value = 146800
$ws.Range("D24").Value = 44800
$ws.Range("E24").Value = 23000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 142500
$ws.Range("E26").Value = 123800
$ws.Range("D27").Value = 142500
$ws.Range("E27").Value = 123800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 9600
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -53400
$ws.Range("E32").Value = -6400
$ws.Range("D33").Value = 152100
$ws.Range("E33").Value = 123800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 152100
$ws.Range("E35").Value = 123800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 904200
$ws.Range("E41").Value = 1398400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 473400
$ws.Range("E43").Value = 378700
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 57500
$ws.Range("E45").Value = 80200
$ws.Range("D46").Value = 1435100
$ws.Range("E46").Value = 1857300
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 90900
$ws.Range("E48").Value = 82700
$ws.Range("D49").Value = 1826600
$ws.Range("E49").Value = 1834400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 35400
$ws.Range("E52").Value = 31900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 3388000
$ws.Range("E54").Value = 3806400
$ws.Range("D57").Value = 3900
$ws.Range("E57").Value = 2200
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 805100
$ws.Range("E59").Value = 672900
$ws.Range("D60").Value = 809000
$ws.Range("E60").Value = 675100
$ws.Range("D61").Value = 2575500
$ws.Range("E61").Value = 2574600
$ws.Range("D62").Value = 169900
$ws.Range("E62").Value = 178300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 3554400
$ws.Range("E66").Value = 3428000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1857000
$ws.Range("E72").Value = 1755700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -166500
$ws.Range("E76").Value = 378400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 152100
$ws.Range("E81").Value = 123800
$ws.Range("D83").Value = 19900
$ws.Range("E83").Value = 19100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 173200
$ws.Range("E89").Value = 143800
$ws.Range("D91").Value = -17200
$ws.Range("E91").Value = -8600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 40000
$ws.Range("E94").Value = -13100
$ws.Range("D96").Value = -50400
$ws.Range("E96").Value = -51700
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -707100
$ws.Range("E100").Value = -97800
$ws.Range("D101").Value = -400
$ws.Range("E101").Value = -2200
$ws.Range("D102").Value = -494200
$ws.Range("E102").Value = 30800
